$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) ---
$ws.Cells.Item(1, 1).Value = "venue"
$ws.Cells.Item(1, 2).Value = "date"
$ws.Cells.Item(1, 3).Value = "result"
$ws.Cells.Item(1, 4).Value = "ownTeam"
$ws.Cells.Item(1, 5).Value = "oppTeam"
$ws.Cells.Item(1, 6).Value = "batsman"
$ws.Cells.Item(1, 7).Value = "totalRuns"
$ws.Cells.Item(1, 8).Value = "totalBalls"
$ws.Cells.Item(1, 9).Value = "total4s"
$ws.Cells.Item(1, 10).Value = "total6s"
$ws.Cells.Item(1, 11).Value = "sr"

# --- Data rows (2-7): text columns (A-F) ---
$ws.Cells.Item(2, 1).Value = " Sharjah"
$ws.Cells.Item(2, 2).Value = " October 23 2020"
$ws.Cells.Item(2, 3).Value = "Mumbai won by 10 wickets (with 46 balls remaining)"
$ws.Cells.Item(2, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(2, 5).Value = "Mumbai Indians"
$ws.Cells.Item(2, 6).Value = "Ruturaj Gaikwad "

$ws.Cells.Item(3, 1).Value = " Abu Dhabi"
$ws.Cells.Item(3, 2).Value = " November 01 2020"
$ws.Cells.Item(3, 3).Value = "Super Kings won by 9 wickets (with 7 balls remaining)"
$ws.Cells.Item(3, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(3, 5).Value = "Kings XI Punjab"
$ws.Cells.Item(3, 6).Value = "Ruturaj Gaikwad "

$ws.Cells.Item(4, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(4, 2).Value = " September 25 2020"
$ws.Cells.Item(4, 3).Value = "Capitals won by 44 runs"
$ws.Cells.Item(4, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(4, 5).Value = "Delhi Capitals"
$ws.Cells.Item(4, 6).Value = "Ruturaj Gaikwad "

$ws.Cells.Item(5, 1).Value = " Sharjah"
$ws.Cells.Item(5, 2).Value = " September 22 2020"
$ws.Cells.Item(5, 3).Value = "Royals won by 16 runs"
$ws.Cells.Item(5, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(5, 5).Value = "Rajasthan Royals"
$ws.Cells.Item(5, 6).Value = "Ruturaj Gaikwad "

$ws.Cells.Item(6, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(6, 2).Value = " October 29 2020"
$ws.Cells.Item(6, 3).Value = "Super Kings won by 6 wickets"
$ws.Cells.Item(6, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(6, 5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(6, 6).Value = "Ruturaj Gaikwad "

$ws.Cells.Item(7, 1).Value = " Dubai (DSC)"
$ws.Cells.Item(7, 2).Value = " October 25 2020"
$ws.Cells.Item(7, 3).Value = "Super Kings won by 8 wickets (with 8 balls remaining)"
$ws.Cells.Item(7, 4).Value = "Chennai Super Kings"
$ws.Cells.Item(7, 5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(7, 6).Value = "Ruturaj Gaikwad "

# --- Numeric-looking columns (G-K) must stay stored as text, matching the source data ---
$numRng = $ws.Range("G2:K7")
$numRng.NumberFormat = "@"
$ws.Cells.Item(2, 7).Value = "0"
$ws.Cells.Item(2, 8).Value = "5"
$ws.Cells.Item(2, 9).Value = "0"
$ws.Cells.Item(2, 10).Value = "0"
$ws.Cells.Item(2, 11).Value = "0.00"
$ws.Cells.Item(3, 7).Value = "62"
$ws.Cells.Item(3, 8).Value = "49"
$ws.Cells.Item(3, 9).Value = "6"
$ws.Cells.Item(3, 10).Value = "1"
$ws.Cells.Item(3, 11).Value = "126.53"
$ws.Cells.Item(4, 7).Value = "5"
$ws.Cells.Item(4, 8).Value = "10"
$ws.Cells.Item(4, 9).Value = "0"
$ws.Cells.Item(4, 10).Value = "0"
$ws.Cells.Item(4, 11).Value = "50.00"
$ws.Cells.Item(5, 7).Value = "0"
$ws.Cells.Item(5, 8).Value = "1"
$ws.Cells.Item(5, 9).Value = "0"
$ws.Cells.Item(5, 10).Value = "0"
$ws.Cells.Item(5, 11).Value = "0.00"
$ws.Cells.Item(6, 7).Value = "72"
$ws.Cells.Item(6, 8).Value = "53"
$ws.Cells.Item(6, 9).Value = "6"
$ws.Cells.Item(6, 10).Value = "2"
$ws.Cells.Item(6, 11).Value = "135.84"
$ws.Cells.Item(7, 7).Value = "65"
$ws.Cells.Item(7, 8).Value = "51"
$ws.Cells.Item(7, 9).Value = "4"
$ws.Cells.Item(7, 10).Value = "3"
$ws.Cells.Item(7, 11).Value = "127.45"
$numRng.Style = "Normal"
